# Aufgabenbereiche.xlsx - add "Priorität" column and a detailed comment for the
# Express-API row (auth work for the node/express API), per:
# "Added auth into nodeexpressapi"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Long comment for the Express-API row (D6) -----------------------------
# Written first so this string lands in sharedStrings.xml before "Priorität".
$longComment = "#Sicherheitsaspekte : Connectstring, Einloggdaten in .env speichern  #Connectiondaten müssen auf aufgesetze Datenbanken angepasst werden, möglich wär auch lokal`n#Datenschemas erstellen"
$ws.Range("D6").Value = $longComment
$ws.Range("D6").WrapText = $true

# Widen column D and grow row 6 so the wrapped comment is readable.
$ws.Columns.Item(4).ColumnWidth = 59.76
$ws.Rows.Item(6).RowHeight = 60

# --- New "Priorität" column (E) ---------------------------------------------
$ws.Range("E1").Value = "Priorität"
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 4
$ws.Range("E8").Value = 3

# --- Simon is also responsible for Login(cookies) (row 10) -----------------
$ws.Range("B10").Value = "Simon"

# Leave the selection on the cell that was last edited in the source file.
$ws.Range("D6").Select() | Out-Null
